# Update multiple stat files - add Austin Hooper, Chris Manhertz, and
# Marcedes Lewis rows to the TE aggregate sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data to append (rows 5-13) ---
# Player, Season Group, Y/R, Y/Tgt, Succ%

# Austin Hooper (yellow highlight, same as new fill 00FFFFBD)
$ws.Cells.Item(5,1).Value = "Austin Hooper"
$ws.Cells.Item(5,2).Value = "Group1"
$ws.Cells.Item(5,3).Value = 9.700000000000001
$ws.Cells.Item(5,4).Value = 6.666666666666667
$ws.Cells.Item(5,5).Value = 56.33333333333334

$ws.Cells.Item(6,1).Value = "Austin Hooper"
$ws.Cells.Item(6,2).Value = "Group2"
$ws.Cells.Item(6,3).Value = 10.26666666666667
$ws.Cells.Item(6,4).Value = 7.666666666666667
$ws.Cells.Item(6,5).Value = 61.33333333333334

$ws.Cells.Item(7,1).Value = "Austin Hooper"
$ws.Cells.Item(7,2).Value = "Difference"
$ws.Cells.Item(7,3).Value = 0.5666666666666664
$ws.Cells.Item(7,4).Value = 1
$ws.Cells.Item(7,5).Value = 5

# Chris Manhertz (reuses existing green highlight, same as Hayden Hurst rows)
$ws.Cells.Item(8,1).Value = "Chris Manhertz"
$ws.Cells.Item(8,2).Value = "Group1"
$ws.Cells.Item(8,3).Value = 10.5
$ws.Cells.Item(8,4).Value = 8.466666666666667
$ws.Cells.Item(8,5).Value = 72.23333333333333

$ws.Cells.Item(9,1).Value = "Chris Manhertz"
$ws.Cells.Item(9,2).Value = "Group2"
$ws.Cells.Item(9,3).Value = 8.333333333333334
$ws.Cells.Item(9,4).Value = 5.600000000000001
$ws.Cells.Item(9,5).Value = 54.16666666666666

$ws.Cells.Item(10,1).Value = "Chris Manhertz"
$ws.Cells.Item(10,2).Value = "Difference"
$ws.Cells.Item(10,3).Value = -2.166666666666666
$ws.Cells.Item(10,4).Value = -2.866666666666666
$ws.Cells.Item(10,5).Value = -18.06666666666667

# Marcedes Lewis (yellow highlight)
$ws.Cells.Item(11,1).Value = "Marcedes Lewis"
$ws.Cells.Item(11,2).Value = "Group1"
$ws.Cells.Item(11,3).Value = 10.13333333333333
$ws.Cells.Item(11,4).Value = 7.366666666666667
$ws.Cells.Item(11,5).Value = 57.93333333333333

$ws.Cells.Item(12,1).Value = "Marcedes Lewis"
$ws.Cells.Item(12,2).Value = "Group2"
$ws.Cells.Item(12,3).Value = 6.766666666666667
$ws.Cells.Item(12,4).Value = 5.399999999999999
$ws.Cells.Item(12,5).Value = 67.13333333333334

$ws.Cells.Item(13,1).Value = "Marcedes Lewis"
$ws.Cells.Item(13,2).Value = "Difference"
$ws.Cells.Item(13,3).Value = -3.366666666666668
$ws.Cells.Item(13,4).Value = -1.966666666666668
$ws.Cells.Item(13,5).Value = 9.20000000000001

# --- Formatting ---
# New pale-yellow fill (00FFFFBD) for the Austin Hooper / Marcedes Lewis blocks
$ws.Range("A5:E7").Interior.Color = 12451839
$ws.Range("A11:E13").Interior.Color = 12451839

# Reuse the existing pale-green highlight (already used for Hayden Hurst, rows 2-4)
# for the Chris Manhertz block, by copying its format.
$ws.Range("A2:E2").Copy()
$ws.Range("A8:E10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
